# Generate Report for Handoff
# Updates the localization-status workbook: statuses flip from
# "In Translation" to "Ready for handoff", the associated timestamps
# advance by about a minute, and the (now-wider) status/date columns on
# each sheet grow to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-13 16:53:26"

# --- zh-cn sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-13 16:53:18"

# --- de-de sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime (date unchanged for de-de)
$dede.Range("C2").Value = "Ready for handoff"

# --- Column widths: the longer "Ready for handoff" text widens the
# status columns (and their mirrored width on the Overview sheet).
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
